$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# 1) Unmerge the two existing merged regions so the row-insert below is
#    unambiguous, then insert a new blank row at row 3. Rows 3-9 shift
#    down to 4-10 (the SUM formula auto-adjusts from F2:F8 to F2:F9,
#    and the sheet dimension grows from B2:F9 to B2:F10).
# ----------------------------------------------------------------------
$ws.Range("D5:D6").UnMerge()
$ws.Range("C5:C7").UnMerge()
$ws.Rows("3:3").Insert()

# ----------------------------------------------------------------------
# 2) The newly inserted row 3 has no useful formatting yet - give it the
#    same look as the rest of the data rows by copying formats only from
#    row 4 (which still carries the original body-row style at this point).
# ----------------------------------------------------------------------
$ws.Range("B4:F4").Copy()
$ws.Range("B3:F3").PasteSpecial(-4122)   # xlPasteFormats

# ----------------------------------------------------------------------
# 3) Write the new "Summary of Changes" row (row 3) content.
# ----------------------------------------------------------------------
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "Summary of Changes"
$ws.Range("D3").Value = "ARF02"
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = 5.5

# ----------------------------------------------------------------------
# 4) Row 4 now holds what used to be row 3's D/E/F values ("---", the
#    "aaaa..." item text, 6.6). Its category cell moved to row 5, so
#    clear the stray "Summary of Changes" text that the insert left
#    behind and fix up the SN column.
# ----------------------------------------------------------------------
$ws.Range("C4").ClearContents()
$ws.Range("B4").Value = 2

# ----------------------------------------------------------------------
# 5) Renumber the SN column for the remaining shifted rows (1..7 total).
# ----------------------------------------------------------------------
$ws.Range("B5").Value = 3
$ws.Range("B6").Value = 4
$ws.Range("B7").Value = 5
$ws.Range("B8").Value = 6
$ws.Range("B9").Value = 7

# ----------------------------------------------------------------------
# 6) Recreate the merges:
#      C3:C4  (new - "Summary of Changes" spans the inserted row)
#      D6:D7  (was D5:D6, shifted down by the insert)
#      C6:C8  (was C5:C7, shifted down by the insert)
# ----------------------------------------------------------------------
$ws.Range("C3:C4").MergeCells = $true
$ws.Range("D6:D7").MergeCells = $true
$ws.Range("C6:C8").MergeCells = $true

# ----------------------------------------------------------------------
# 7) Merging re-derives a style for every cell it touches, which clobbers
#    the plain "top / middle / bottom of a bordered block" look the rest
#    of the sheet uses. Restore it:
#       - "top" cells (have content) -> same look as any other text cell
#       - "bottom" cells (blank)     -> thin left+right+bottom border
#       - "middle" cell  (blank)     -> thin left+right border only
# ----------------------------------------------------------------------
$ws.Range("B3").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("B6").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("B6").Copy()
$ws.Range("C6").PasteSpecial(-4122)

foreach ($addr in @("C4", "D7", "C8")) {
    $ws.Range("A1").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($addr).Borders.Item(7).LineStyle = 1    # left
    $ws.Range($addr).Borders.Item(10).LineStyle = 1   # right
    $ws.Range($addr).Borders.Item(9).LineStyle = 1    # bottom
}

$ws.Range("A1").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("C7").Borders.Item(7).LineStyle = 1    # left
$ws.Range("C7").Borders.Item(10).LineStyle = 1   # right

# re-assert content on the merge anchor cells in case PasteSpecial touched it
$ws.Range("C3").Value = "Summary of Changes"
$ws.Range("D6").Value = "ARF02"
$ws.Range("C6").Value = "Interface Impact"

Write-Output "edit complete"
